$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header row (row 1) with new uppercase/underscore style names
$ws.Range("A1").Value = "Sample_ID"
$ws.Range("B1").Value = "Toluene_Area"
$ws.Range("C1").Value = "Toluene_Amt"
$ws.Range("D1").Value = "Sample_Vol"
$ws.Range("E1").Value = "No_Hours"
$ws.Range("F1").Value = "No_Flowers"
$ws.Range("G1").Value = "Fresh_Mass"
$ws.Range("H1").Value = "Dry_Mass"

# Update sample id values in column A (rows 2-6) to use "Samp." prefix
$ws.Range("A2").Value = "Samp.rm3-1"
$ws.Range("A3").Value = "Samp.rm3-2"
$ws.Range("A4").Value = "Samp.rm3-4-1"
$ws.Range("A5").Value = "Samp.rm3-7"
$ws.Range("A6").Value = "Samp.rm3-8"

# Update the active selection on the sheet
$ws.Range("A10").Select()
